# Update Name of Algo
# Update column B values (KNN imputation results) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    8   = 6.185
    10  = 5.825000000000001
    12  = 5.315
    18  = 5.125999999999999
    25  = 6.396
    37  = 8.643000000000001
    55  = 4.684
    68  = 5.356
    77  = 5.737
    78  = 7.571000000000001
    79  = 5.542
    80  = 7.845999999999999
    81  = 6.452000000000001
    82  = 5.659000000000001
    84  = 5.427000000000001
    101 = 6.947
    102 = 7.316
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
